$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" '71.890.50'
Set-TextValue "E2" '  -0.76%  '
Set-TextValue "D3" '3.998.63'
Set-TextValue "E3" '  -1.23%  '
Set-TextValue "D4" '0.997'
Set-TextValue "E4" '  -0.31%  '
Set-TextValue "D5" '545.74'
Set-TextValue "E5" '  +4.73%  '
Set-TextValue "D6" '150.15'
Set-TextValue "E6" '  +1.75%  '
Set-TextValue "D7" '0.703'
Set-TextValue "E7" '  +12.73%  '
Set-TextValue "D8" '0.999'
Set-TextValue "E8" '  -0.07%  '
Set-TextValue "D9" '0.746'
Set-TextValue "E9" '  +1.16%  '
Set-TextValue "D10" '0.170'
Set-TextValue "E10" '  -2.95%  '
Set-TextValue "B11" 'Avalanche'
Set-TextValue "C11" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D11" '51.76'
Set-TextValue "E11" '  +8.56%  '
Set-TextValue "B12" 'ShibaInu'
Set-TextValue "C12" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D12" '0.0000324'
Set-TextValue "E12" '  -3.19%  '
Set-TextValue "D13" '10.70'
Set-TextValue "E13" '  -1.87%  '
Set-TextValue "D14" '4.612.67'
Set-TextValue "E14" '  -1.39%  '
Set-TextValue "D15" '3.980.79'
Set-TextValue "E15" '  -1.49%  '
Set-TextValue "D16" '14.07'
Set-TextValue "E16" '  -0.77%  '
Set-TextValue "D17" '20.47'
Set-TextValue "E17" '  -3.51%  '
Set-TextValue "E18" '  -0.26%  '
Set-TextValue "E19" '  -2.21%  '
Set-TextValue "D20" '71.643.98'
Set-TextValue "E20" '  -0.99%  '
Set-TextValue "D21" '429.62'
Set-TextValue "E21" '  -1.80%  '
Set-TextValue "D22" '97.24'
Set-TextValue "E22" '  -0.71%  '
Set-TextValue "D23" '3.51'
Set-TextValue "E23" '  -0.75%  '
Set-TextValue "D24" '4.26'
Set-TextValue "E24" '  +5.87%  '
Set-TextValue "D25" '14.32'
Set-TextValue "E25" '  -3.31%  '
Set-TextValue "D26" '11.44'
Set-TextValue "E26" '  -4.29%  '
Set-TextValue "D27" '10.73'
Set-TextValue "E27" '  -5.01%  '
Set-TextValue "D28" '5.84'
Set-TextValue "E28" '  +0.97%  '
Set-TextValue "B29" 'Toncoin'
Set-TextValue "C29" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D29" '3.64'
Set-TextValue "E29" '  +17.64%  '
Set-TextValue "B30" 'EthereumClassic'
Set-TextValue "C30" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D30" '36.73'
Set-TextValue "E30" '  -1.09%  '
Set-TextValue "B31" 'Cosmos'
Set-TextValue "C31" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D31" '13.40'
Set-TextValue "E31" '  -1.04%  '
Set-TextValue "B32" 'InjectiveProtocol'
Set-TextValue "C32" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D32" '49.00'
Set-TextValue "E32" '  +20.38%  '
Set-TextValue "D33" '7.24'
Set-TextValue "E33" '  +2.90%  '
Set-TextValue "B34" 'Hedera'
Set-TextValue "C34" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D34" '0.130'
Set-TextValue "E34" '  +1.18%  '
Set-TextValue "D35" '676.94'
Set-TextValue "E35" '  -2.33%  '
Set-TextValue "D36" '65.66'
Set-TextValue "E36" '  -3.72%  '
Set-TextValue "D37" '0.441'
Set-TextValue "E37" '  +0.56%  '
Set-TextValue "B38" 'PEPE'
Set-TextValue "C38" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D38" '0.0₃0836'
Set-TextValue "E38" '  -6.49%  '
Set-TextValue "B39" 'Kaspa'
Set-TextValue "C39" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D39" '0.151'
Set-TextValue "E39" '  -1.93%  '
Set-TextValue "D40" '3.41'
Set-TextValue "E40" '  -6.59%  '
Set-TextValue "E41" '  -0.02%  '
Set-TextValue "D42" '3.32'
Set-TextValue "E42" '  +5.07%  '
Set-TextValue "E43" '  +0.28%  '
Set-TextValue "E44" '  -0.59%  '
Set-TextValue "D45" '0.150'
Set-TextValue "E45" '  +2.82%  '
Set-TextValue "B46" 'Fetch.AI'
Set-TextValue "C46" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D46" '2.72'
Set-TextValue "E46" '  -2.09%  '
Set-TextValue "B47" 'THORChain'
Set-TextValue "C47" 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue "D47" '9.83'
Set-TextValue "E47" '  +8.51%  '
Set-TextValue "B48" 'FLOKI'
Set-TextValue "C48" 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue "D48" '0.000282'
Set-TextValue "E48" '  +3.24%  '
Set-TextValue "B49" 'ApeXProtocol'
Set-TextValue "C49" 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue "D49" '3.31'
Set-TextValue "E49" '  -5.28%  '
Set-TextValue "B50" 'Stacks'
Set-TextValue "C50" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D50" '3.00'
Set-TextValue "E50" '  -3.62%  '
Set-TextValue "D51" '144.59'
Set-TextValue "E51" '  +1.47%  '
